$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 632: new-cases count revised down by 1 ---
$ws.Range("C632").Value = 183

# --- Row 633: new-cases count revised down by 1 ---
$ws.Range("C633").Value = 169

# --- Row 634: C634 itself is unchanged (101); only the downstream
#     cumulative formula in B634 shifts because B633 changed. ---

# --- Row 635 ---
$ws.Range("C635").Value = 100

# --- Row 636 ---
$ws.Range("C636").Value = 302

# --- Row 637 ---
$ws.Range("C637").Value = 247

# --- Row 638 ---
$ws.Range("C638").Value = 269
$ws.Range("E638").Value = 8
$ws.Range("G638").Value = 51

# --- Row 639 (fills in previously-missing ICU/intubated/hospital-death data) ---
$ws.Range("C639").Value = 271
$ws.Range("E639").Value = 10
$ws.Range("F639").Value = 4
$ws.Range("G639").Value = 49

# --- Row 640 (was a blank placeholder row; now populated with real data) ---
$ws.Range("C640").Value = 263
$ws.Range("E640").Value = 11
$ws.Range("F640").Value = 3
$ws.Range("G640").Value = 48

# --- Row 641 (was a blank placeholder row; now populated with real data) ---
$ws.Range("C641").Value = 146
$ws.Range("E641").Value = 11
$ws.Range("F641").Value = 3
$ws.Range("G641").Value = 52

# --- Row 642 (was a blank placeholder row; now populated with real data) ---
$ws.Range("C642").Value = 85
$ws.Range("E642").Value = 12
$ws.Range("F642").Value = 4
$ws.Range("G642").Value = 56

# --- Row 643 (was a blank placeholder row; now populated with real data) ---
$ws.Range("C643").Value = 10
$ws.Range("E643").Value = 10
$ws.Range("F643").Value = 1
$ws.Range("G643").Value = 53

# --- Columns L (deaths at hospital) & M (deaths extra-hospital) are
#     formatted as Text (numFmtId 49) but must keep storing genuine
#     numbers, matching the rest of the column. Writing straight into a
#     Text-formatted cell stores the digits as a string, so flip the
#     cell to General, write the number, then flip the format back to
#     Text so the stored <v> stays numeric while the visible format is
#     unchanged. ---
$lOnes = $ws.Range("L635")
$lOnes.NumberFormat = "General"
$lOnes.Value = 1
$lOnes.NumberFormat = "@"

$lZeros = $ws.Range("L640:L643")
$lZeros.NumberFormat = "General"
$lZeros.Value = 0
$lZeros.NumberFormat = "@"

$mOnes = $ws.Range("M639")
$mOnes.NumberFormat = "General"
$mOnes.Value = 1
$mOnes.NumberFormat = "@"

$mZeros = $ws.Range("M640:M643")
$mZeros.NumberFormat = "General"
$mZeros.Value = 0
$mZeros.NumberFormat = "@"

# --- Restore the active-cell selection saved with the workbook (A2) ---
$ws.Range("A2").Select()
